$wb = $excel.ActiveWorkbook

# --- Release mCSD 3.9.0 with CP integrated ---

# Rename the two "Include from ..." include-sheets to their generic
# "Include #N" titles (IG Publisher renders tab names this way once the
# content no longer fits the 31-char sheet-name limit nicely).
$wb.Worksheets.Item("Include from mCSD Organizatio").Name = "Include #0"
$wb.Worksheets.Item("Include from Organization Aff").Name = "Include #1"

# Update the Metadata sheet with the refreshed ValueSet metadata.
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "3.9.0"
$ws.Range("B7").Value = "false"
$ws.Range("B8").Value = "2024-12-02T17:05:26-06:00"
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"
$ws.Range("B13").Value = "Global (Whole world)"
